# Refresh the crypto price (D) / 1h volume change (E) figures for each
# ranking row, matching the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.646.89"
$ws.Range("E2").Value = "  -4.91%  "
$ws.Range("D3").Value = "2.207.01"
$ws.Range("E3").Value = "  -5.85%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'244.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.97%  "
$ws.Range("D7").Value = "'70.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.26%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -8.43%  "
$ws.Range("D10").Value = "'37.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.00%  "
$ws.Range("E11").Value = "  -6.12%  "
$ws.Range("D12").Value = "'57.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.99%  "
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "'6.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.43%  "
$ws.Range("D15").Value = "2.535.17"
$ws.Range("E15").Value = "  -5.93%  "
$ws.Range("D16").Value = "'14.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.59%  "
$ws.Range("D17").Value = "'0.837"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.16%  "
$ws.Range("D18").Value = "2.207.32"
$ws.Range("E18").Value = "  -5.52%  "
$ws.Range("D19").Value = "41.567.21"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -6.99%  "
$ws.Range("D21").Value = "'73.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.25%  "
$ws.Range("D22").Value = "'6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.74%  "
$ws.Range("D23").Value = "'233.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.28%  "
$ws.Range("E24").Value = "  +10.54%  "
$ws.Range("E26").Value = "  -5.08%  "
$ws.Range("D27").Value = "'2.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "'9.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("D30").Value = "'168.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").Value = "'20.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.33%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.72%  "
$ws.Range("E33").Value = "  -6.96%  "
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "'5.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("D36").Value = "'4.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.60%  "
$ws.Range("D37").Value = "'3.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").Value = "'23.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.05%  "
$ws.Range("E39").Value = "  -4.93%  "
$ws.Range("D40").Value = "'0.0270"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'5.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.81%  "
$ws.Range("D42").Value = "'64.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "'8.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("E44").Value = "  -10.65%  "
$ws.Range("D45").Value = "'0.191"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("D46").Value = "'0.0993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'4.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.32%  "
$ws.Range("E49").Value = "  +8.08%  "
$ws.Range("E50").Value = "  -4.15%  "
$ws.Range("D51").Value = "0.0₃0147"
$ws.Range("E51").Value = "  +12.77%  "
